$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 495.70587
$ws.Range("I9").Value = 456.2
$ws.Range("J9").Value = 552.1429000000001
$ws.Range("K9").Value = 456.2
$ws.Range("L9").Value = 552.1429000000001
$ws.Range("M9").Value = -287.2
$ws.Range("N9").Value = -890.1429000000001
$ws.Range("H11").Value = 1860.2
$ws.Range("I11").Value = 1860.2
$ws.Range("K11").Value = 1860.2
$ws.Range("M11").Value = -1720.2
$ws.Range("H51").Value = 5038.5557
$ws.Range("I51").Value = 6500
$ws.Range("J51").Value = 4855.875
$ws.Range("K51").Value = 6500
$ws.Range("L51").Value = 4855.875
$ws.Range("M51").Value = -6016
$ws.Range("N51").Value = -5823.875
$ws.Range("H70").Value = 116670270
$ws.Range("I70").Value = 125000500
$ws.Range("J70").Value = 111116780
$ws.Range("K70").Value = 375001500
$ws.Range("L70").Value = 333350340
$ws.Range("M70").Value = -375001230
$ws.Range("N70").Value = -333350880
$ws.Range("H73").Value = 116670270
$ws.Range("I73").Value = 125000500
$ws.Range("J73").Value = 111116780
$ws.Range("K73").Value = 375001500
$ws.Range("L73").Value = 333350340
$ws.Range("M73").Value = -375000564
$ws.Range("N73").Value = -333352212
$ws.Range("H80").Value = 33993.734
$ws.Range("I80").Value = 11489.223
$ws.Range("J80").Value = 67750.5
$ws.Range("K80").Value = 34467.669
$ws.Range("L80").Value = 203251.5
$ws.Range("M80").Value = -33469.669
$ws.Range("N80").Value = -205247.5
$ws.Range("H83").Value = 33993.734
$ws.Range("I83").Value = 11489.223
$ws.Range("J83").Value = 67750.5
$ws.Range("K83").Value = 103403.007
$ws.Range("L83").Value = 609754.5
$ws.Range("M83").Value = -98411.007
$ws.Range("N83").Value = -619738.5
$ws.Range("H86").Value = 65587824
$ws.Range("I86").Value = 112500590
$ws.Range("K86").Value = 112500590
$ws.Range("M86").Value = -112499467
$ws.Range("H89").Value = 65587824
$ws.Range("I89").Value = 112500590
$ws.Range("K89").Value = 562502950
$ws.Range("M89").Value = -562497334
$ws.Range("H101").Value = 289.3846
$ws.Range("I101").Value = 287.7143
$ws.Range("J101").Value = 291.33334
$ws.Range("K101").Value = 863.1428999999999
$ws.Range("L101").Value = 874.0000200000001
$ws.Range("M101").Value = 758.8571000000001
$ws.Range("N101").Value = -4118.00002
$ws.Range("H118").Value = 494.16666
$ws.Range("I118").Value = 503
$ws.Range("K118").Value = 1509
$ws.Range("M118").Value = 148
$ws.Range("H137").Value = 7178.7
$ws.Range("I137").Value = 3060.3845
$ws.Range("K137").Value = 9181.1535
$ws.Range("M137").Value = -6631.1535
$ws.Range("H138").Value = 2004838.9
$ws.Range("J138").Value = 3339975.5
$ws.Range("L138").Value = 10019926.5
$ws.Range("N138").Value = -10030206.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3233122.8
$ws.Range("I32").Value = 3514098.5
$ws.Range("K32").Value = 3514098.5
$ws.Range("M32").Value = -3513811.5
$ws.Range("H61").Value = 55558804
$ws.Range("I61").Value = 2536.4285
$ws.Range("K61").Value = 2536.4285
$ws.Range("M61").Value = -2324.4285
$ws.Range("H133").Value = 116841
$ws.Range("J133").Value = 116841
$ws.Range("L133").Value = 116841
$ws.Range("N133").Value = -121901
$ws.Range("H136").Value = 55558804
$ws.Range("I136").Value = 2536.4285
$ws.Range("K136").Value = 7609.2855
$ws.Range("M136").Value = -5059.2855
$ws.Range("H140").Value = 57716
$ws.Range("J140").Value = 57716
$ws.Range("L140").Value = 57716
$ws.Range("N140").Value = -68076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 75000
$ws.Range("J74").Value = 75000
$ws.Range("L74").Value = 75000
$ws.Range("N74").Value = -76872
$ws.Range("H77").Value = 75000
$ws.Range("J77").Value = 75000
$ws.Range("L77").Value = 225000
$ws.Range("N77").Value = -234360
$ws.Range("H81").Value = 90593.336
$ws.Range("J81").Value = 90593.336
$ws.Range("L81").Value = 90593.336
$ws.Range("N81").Value = -92715.336
$ws.Range("H84").Value = 90593.336
$ws.Range("J84").Value = 90593.336
$ws.Range("L84").Value = 271780.008
$ws.Range("N84").Value = -282388.008
$ws.Range("H105").Value = 3306.95
$ws.Range("I105").Value = 2626.8462
$ws.Range("J105").Value = 4570
$ws.Range("K105").Value = 2626.8462
$ws.Range("L105").Value = 4570
$ws.Range("M105").Value = -879.8462
$ws.Range("N105").Value = -8064
$ws.Range("H134").Value = 3679824.5
$ws.Range("I134").Value = 4466668
$ws.Range("J134").Value = 7887
$ws.Range("K134").Value = 13400004
$ws.Range("L134").Value = 23661
$ws.Range("M134").Value = -13397469
$ws.Range("N134").Value = -28731
$ws.Range("H139").Value = 61333
$ws.Range("J139").Value = 67599.60000000001
$ws.Range("L139").Value = 67599.60000000001
$ws.Range("N139").Value = -77879.60000000001
$ws.Range("H141").Value = 76994.25
$ws.Range("J141").Value = 76994.25
$ws.Range("L141").Value = 76994.25
$ws.Range("N141").Value = -87354.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 4544.481
$ws.Range("I31").Value = 1825.5385
$ws.Range("J31").Value = 7263.423
$ws.Range("K31").Value = 1825.5385
$ws.Range("L31").Value = 7263.423
$ws.Range("M31").Value = -1530.5385
$ws.Range("N31").Value = -7853.423
$ws.Range("H34").Value = 4544.481
$ws.Range("I34").Value = 1825.5385
$ws.Range("J34").Value = 7263.423
$ws.Range("K34").Value = 1825.5385
$ws.Range("L34").Value = 7263.423
$ws.Range("M34").Value = -1623.5385
$ws.Range("N34").Value = -7667.423
$ws.Range("H132").Value = 5358.88
$ws.Range("I132").Value = 3690.923
$ws.Range("J132").Value = 7165.8335
$ws.Range("K132").Value = 11072.769
$ws.Range("L132").Value = 21497.5005
$ws.Range("M132").Value = -8542.769
$ws.Range("N132").Value = -26557.5005
$ws.Range("H134").Value = 3719.3777
$ws.Range("I134").Value = 1629.6316
$ws.Range("K134").Value = 4888.8948
$ws.Range("M134").Value = -2353.8948
$ws.Range("H135").Value = 64500
$ws.Range("I135").Value = 49000
$ws.Range("J135").Value = 80000
$ws.Range("K135").Value = 49000
$ws.Range("L135").Value = 80000
$ws.Range("M135").Value = -43930
$ws.Range("N135").Value = -90140
$ws.Range("H137").Value = 74999
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1030.1111
$ws.Range("I5").Value = 783.25
$ws.Range("J5").Value = 3005
$ws.Range("K5").Value = 2349.75
$ws.Range("L5").Value = 9015
$ws.Range("M5").Value = -2237.75
$ws.Range("N5").Value = -9239
$ws.Range("H113").Value = 5298.3076
$ws.Range("I113").Value = 683.3333
$ws.Range("J113").Value = 6682.8
$ws.Range("K113").Value = 2049.9999
$ws.Range("L113").Value = 20048.4
$ws.Range("M113").Value = 120.0001000000002
$ws.Range("N113").Value = -24388.4
$ws.Range("H135").Value = 1030.1111
$ws.Range("I135").Value = 783.25
$ws.Range("J135").Value = 3005
$ws.Range("K135").Value = 7049.25
$ws.Range("L135").Value = 27045
$ws.Range("M135").Value = -4514.25
$ws.Range("N135").Value = -32115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 79215.234
$ws.Range("I80").Value = 2334
$ws.Range("K80").Value = 2334
$ws.Range("M80").Value = -1336
$ws.Range("H83").Value = 79215.234
$ws.Range("I83").Value = 2334
$ws.Range("K83").Value = 11670
$ws.Range("M83").Value = -6678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9943.564
$ws.Range("I136").Value = 2773.913
$ws.Range("J136").Value = 20249.938
$ws.Range("K136").Value = 8321.739
$ws.Range("L136").Value = 60749.814
$ws.Range("M136").Value = -5771.739
$ws.Range("N136").Value = -65849.814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9645.875
$ws.Range("I132").Value = 14714
$ws.Range("K132").Value = 44142
$ws.Range("M132").Value = -41612
$ws.Range("H136").Value = 24638066
$ws.Range("I136").Value = 43479440
$ws.Range("J136").Value = 562978.1
$ws.Range("K136").Value = 130438320
$ws.Range("L136").Value = 1688934.3
$ws.Range("M136").Value = -130435770
$ws.Range("N136").Value = -1694034.3
$ws.Range("H140").Value = 72500
$ws.Range("J140").Value = 72500
$ws.Range("L140").Value = 72500
$ws.Range("N140").Value = -82860
